$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Title paragraph: split the single run of text into one run per
#    word/space, still as plain text runs with no character formatting.
# ---------------------------------------------------------------------------
$titlePara = $d.Paragraphs.Item(1)
$titleXml = '<?xml version="1.0" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="Title"/></w:pPr><w:r><w:t xml:space="preserve">Week</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">7</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">Reading</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">Guide</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">Part</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">1:</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">Confidence</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">Intervals</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$titlePara.Range.InsertXML($titleXml)

# ---------------------------------------------------------------------------
# 2. Tip picture: fix the embedded alt-text / description path baked into
#    the picture's non-visual properties (it pointed at a dev machine's
#    RStudio install; trim it down to the generic quarto path).
# ---------------------------------------------------------------------------
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.InlineShapes.Count -gt 0) {
        $tipXml = '<?xml version="1.0" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:wp="http://schemas.openxmlformats.org/drawingml/2006/wordprocessingDrawing" xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main" xmlns:pic="http://schemas.openxmlformats.org/drawingml/2006/picture" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships"><w:body><w:p><w:pPr><w:pStyle w:val="BodyText"/></w:pPr><w:pPr><w:spacing w:before="0" w:after="0"/><w:textAlignment w:val="center"/></w:pPr><w:r><w:drawing><wp:inline><wp:extent cx="152400" cy="152400"/><wp:effectExtent b="0" l="0" r="0" t="0"/><wp:docPr descr="" title="" id="23" name="Picture"/><a:graphic><a:graphicData uri="http://schemas.openxmlformats.org/drawingml/2006/picture"><pic:pic><pic:nvPicPr><pic:cNvPr descr="/Applications/quarto/share/formats/docx/tip.png" id="24" name="Picture"/><pic:cNvPicPr><a:picLocks noChangeArrowheads="1" noChangeAspect="1"/></pic:cNvPicPr></pic:nvPicPr><pic:blipFill><a:blip r:embed="rId22"/><a:stretch><a:fillRect/></a:stretch></pic:blipFill><pic:spPr bwMode="auto"><a:xfrm><a:off x="0" y="0"/><a:ext cx="152400" cy="152400"/></a:xfrm><a:prstGeom prst="rect"><a:avLst/></a:prstGeom><a:noFill/><a:ln w="9525"><a:noFill/><a:headEnd/><a:tailEnd/></a:ln></pic:spPr></pic:pic></a:graphicData></a:graphic></wp:inline></w:drawing></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">Tip</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
        $p.Range.InsertXML($tipXml)
        break
    }
}

# ---------------------------------------------------------------------------
# 3. "Subtitle" style: re-base it off Normal instead of Title, and pin an
#    explicit (grey) text color so it keeps looking the same now that it no
#    longer inherits Title's color.
# ---------------------------------------------------------------------------
$subtitle = $d.Styles.Item("Subtitle")
$subtitle.BaseStyle = "Normal"
$subtitle.Font.TextColor.ObjectThemeColor = 13  # wdThemeColorText1 -> text1

# ---------------------------------------------------------------------------
# 4. "AbstractTitle" style: pin an explicit blue text color.
# ---------------------------------------------------------------------------
$abstractTitle = $d.Styles.Item("AbstractTitle")
$abstractTitle.Font.Color = 9067060  # wdColor BGR for 345A8A
